$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (columns L, M, N re-labelled) ---
$ws.Range("L1").Value = "MT brut (Rappel)"
$ws.Range("M1").Value = "Taxe (Rappel)"
$ws.Range("N1").Value = "Caution"

# --- Row 2 ---
$ws.Range("A2").Value = "990/PV 01"
$ws.Range("B2").Value = "Point de vente"
$ws.Range("C2").Value = "L3578354"
$ws.Range("D2").Value = "NABIL KAMAL"
$ws.Range("E2").Value = "non"
$ws.Range("F2").Value = "trimestrielle"
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = "--"
$ws.Range("I2").Value = 10000
$ws.Range("J2").Value = "--"
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 10000
$ws.Range("O2").Value = 9000

# --- Row 3 ---
$ws.Range("A3").Value = "044/LF/FES VILLE /AV1"
$ws.Range("B3").Value = "Logement de fonction"
$ws.Range("C3").Value = "K5443645"
$ws.Range("D3").Value = "KHADIJA LALA"
$ws.Range("E3").Value = "non"
$ws.Range("F3").Value = "mensuelle"
$ws.Range("G3").Value = 15
$ws.Range("H3").Value = 10000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = "--"
$ws.Range("O3").Value = 8500

# --- Row 4 ---
$ws.Range("A4").Value = "554/SUP FES 1"
$ws.Range("B4").Value = "Supervision"
$ws.Range("C4").Value = "D524564"
$ws.Range("D4").Value = "SAMIRA TATA"
$ws.Range("E4").Value = "non"
$ws.Range("F4").Value = "mensuelle"
$ws.Range("G4").Value = 15
$ws.Range("H4").Value = 10000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = "--"
$ws.Range("O4").Value = 8500

# --- Row 5 ---
$ws.Range("A5").Value = "800/PV FES 1"
$ws.Range("B5").Value = "Point de vente"
$ws.Range("C5").Value = "P5874857"
$ws.Range("D5").Value = "KARIM JALAL"
$ws.Range("E5").Value = "non"
$ws.Range("F5").Value = "mensuelle"
$ws.Range("G5").Value = 15
$ws.Range("H5").Value = 10000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = "--"
$ws.Range("O5").Value = 8500

# --- Row 6 ---
$ws.Range("A6").Value = " "
$ws.Range("B6").Value = " "
$ws.Range("C6").Value = " "
$ws.Range("D6").Value = " "
$ws.Range("E6").Value = " "
$ws.Range("F6").Value = " "
$ws.Range("G6").Value = " "
$ws.Range("H6").Value = 30000
$ws.Range("I6").Value = 10000
$ws.Range("J6").Value = 4500
$ws.Range("K6").Value = 1000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 34500

# --- Remove the now-unused trailing rows 7 and 8 ---
$ws.Rows("7:8").Delete()
